# Weekly driver report update for 2025-04-29
# Rebuilds the "Good Drivers" table (rows 11-13) with the new column
# layout (good/critical/warning/total sums, adapter/driver split, etc.),
# drops the old bordered/bold header styling + number formatting on that
# block, resizes a handful of columns, and lets the trailing blank rows
# (14-18) fall away naturally once they hold no content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the old "Good Drivers" block entirely (values + styles) so the
# rebuilt cells carry no formatting at all, matching the plain (unstyled)
# cells in the updated report.
$ws.Range("A11:J18").Clear()

# --- Column width tweaks -------------------------------------------------
# Raw stored width = ColumnWidth + 0.83 in this engine's units, so offset
# the target widths accordingly to land on exact integer stored widths.
$ws.Columns.Item(2).ColumnWidth = 14 - 0.83    # B: 15 -> 14
$ws.Columns.Item(5).ColumnWidth = 14 - 0.83    # E: 16 -> 14
$ws.Columns.Item(6).ColumnWidth = 11 - 0.83    # F: 2  -> 11
$ws.Columns.Item(7).ColumnWidth = 23 - 0.83    # G: 2  -> 23
$ws.Columns.Item(8).ColumnWidth = 13 - 0.83    # H: 2  -> 13
$ws.Columns.Item(9).ColumnWidth = 30 - 0.83    # I: 2  -> 30
$ws.Columns.Item(10).ColumnWidth = 16 - 0.83   # J: 2  -> 16

# --- Row 11: new headers --------------------------------------------------
$ws.Range("A11").Value = "adapter-driver"
$ws.Range("B11").Value = "good sum"
$ws.Range("C11").Value = "critical sum"
$ws.Range("D11").Value = "warning sum"
$ws.Range("E11").Value = "client count"
$ws.Range("F11").Value = "total sum"
$ws.Range("G11").Value = "adapter"
$ws.Range("H11").Value = "driver"
$ws.Range("I11").Value = "good roaming calculation (%)"
$ws.Range("J11").Value = "driver vintage"

# --- Row 12: RZ616 Wi-Fi 6E 160MHz - 23.32.2.560 --------------------------
$ws.Range("A12").Value = "RZ616 Wi-Fi 6E 160MHz - 23.32.2.560"
$ws.Range("B12").Value = 13763
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 69
$ws.Range("F12").Value = 13764
$ws.Range("G12").Value = "rz616 wi-fi 6e 160mhz"
$ws.Range("H12").Value = "23.32.2.560"
$ws.Range("I12").Value = 100

# --- Row 13: RZ616 Wi-Fi 6E 160MHz - 3.3.0.908 ----------------------------
$ws.Range("A13").Value = "RZ616 Wi-Fi 6E 160MHz - 3.3.0.908"
$ws.Range("B13").Value = 56907
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 47
$ws.Range("F13").Value = 56907
$ws.Range("G13").Value = "rz616 wi-fi 6e 160mhz"
$ws.Range("H13").Value = "3.3.0.908"
$ws.Range("I13").Value = 100

# --- Driver Vintage (J12/J13): these look like ISO dates, so a plain
# .Value assignment would get silently reinterpreted as a date serial +
# date number format. Route them through a self-quoting text formula and
# then Paste-Values-only over themselves, which collapses them back down
# to plain literal text with no residual formula or number formatting.
$ws.Range("J12").Formula = "=""2023-07-16"""
$ws.Range("J13").Formula = "=""2023-11-26"""
$ws.Range("J12:J13").Copy()
$ws.Range("J12:J13").PasteSpecial(-4163)
